$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for rows 2, 3, 4 in the affected columns
# (A, B, D, E, F, G, H, Q, R). These values get cyclically shifted: the data
# that currently lives in row 3 moves to row 2, row 4's data moves to row 3,
# and row 2's data moves to row 4.

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

$row2 = @{}
$row3 = @{}
$row4 = @{}

foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value2
    $row3[$col] = $ws.Range("${col}3").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}2").Value = $row3[$col]
    $ws.Range("${col}3").Value = $row4[$col]
    $ws.Range("${col}4").Value = $row2[$col]
}
